$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing C8 score value (2.5 -> 2.75) ---
$ws.Cells.Item(8, 3).Value2 = 2.75

# --- Update the tasting-note text shared by D8 & the new D9 (7 C -> 11 C) ---
$secondNote = "Cooled and served @ 11 C. Opening the swing-top gave a puff. Poured clear with a little yeast from conditioning. No foam, light carbonation. A light sour taste, boozy, some light esters and dry mouthfeel."
$ws.Cells.Item(8, 4).Value2 = $secondNote

# --- Row 9: new tasting entry dated 2020-10-13 (serial 44117) ---
$ws.Cells.Item(9, 1).Value2 = 44117
$ws.Cells.Item(9, 2).Formula = "=A9-`$A`$6"
$ws.Cells.Item(9, 3).Value2 = 2.75
$ws.Cells.Item(9, 4).Value2 = $secondNote
$ws.Rows.Item(9).RowHeight = 60

# Re-assert the plain-number format on the day-count formula cell (Excel's date-subtraction
# auto-format would otherwise stamp it with the date number format inherited from column A).
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 10: new tasting entry dated 2020-10-24 (serial 44128) ---
$thirdNote = "Cooled and served @ 11 C. Opening the swing-top gave a puff. Poured clear with a little yeast from conditioning. No foam, light carbonation. A light sour taste, boozy and dry mouthfeel."
$ws.Cells.Item(10, 1).Value2 = 44128
$ws.Cells.Item(10, 2).Formula = "=A10-`$A`$6"
$ws.Cells.Item(10, 3).Value2 = 2.75
$ws.Cells.Item(10, 4).Value2 = $thirdNote
$ws.Rows.Item(10).RowHeight = 45

$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Extend the table with a new blank row 14 (same formatting as row 13) ---
$ws.Range("B13:D13").Copy()
$ws.Range("B14:D14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the active selection to D11 ---
$ws.Range("D11").Select()
